# Commit: "results updated for without stopword"
# Adds a new "Random Forest" result row (row 19) to Sheet 1, continuing the
# "Class Weights / Without Stopwords / SVD + MinMaxScaler + SMOTE 0.3" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# New row 19 values (mirrors the layout of row 18, same Type/Version/SMOTE
# grouping, new Model = "Random Forest").
$ws.Range("A19").Value = "Class Weights"
$ws.Range("B19").Value = "Without Stopwords"
$ws.Range("C19").Value = "SVD + MinMaxScaler + SMOTE 0.3"
$ws.Range("D19").Value = "Random Forest"
$ws.Range("E19").Value = 0.8669
$ws.Range("F19").Value = 0.44
$ws.Range("G19").Value = 0.03
$ws.Range("H19").Value = 0.05
$ws.Range("I19").Value = 0.6583

# Match the formatting of the row above (borders, alignment, number formats).
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(19).RowHeight = 45

# Re-apply values in case paste-special touched anything.
$ws.Range("A19").Value = "Class Weights"
$ws.Range("B19").Value = "Without Stopwords"
$ws.Range("C19").Value = "SVD + MinMaxScaler + SMOTE 0.3"
$ws.Range("D19").Value = "Random Forest"
$ws.Range("E19").Value = 0.8669
$ws.Range("F19").Value = 0.44
$ws.Range("G19").Value = 0.03
$ws.Range("H19").Value = 0.05
$ws.Range("I19").Value = 0.6583

# Leave selection on the newly added cell, as in the saved workbook.
$ws.Range("I19").Select()
